$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.258.13"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "3.838.90"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("E4").Value = "  -0.28%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "448.58"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +7.11%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "147.37"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +15.37%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.626"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +4.85%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.14%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.745"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +4.18%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.158"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.92%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0000326"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -3.87%  "
$ws.Range("E12").Value = "  +10.55%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "10.41"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +4.26%  "
$ws.Range("D14").Value = "4.443.06"
$ws.Range("E14").Value = "  +1.04%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "15.11"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -4.51%  "
$ws.Range("D16").Value = "3.811.13"
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("E17").Value = "  -0.26%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "20.05"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +3.89%  "
$ws.Range("E19").Value = "  +8.44%  "
$ws.Range("D20").Value = "67.272.74"
$ws.Range("E20").Value = "  +0.92%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "426.94"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +6.47%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "14.72"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +4.34%  "
$ws.Range("E23").Value = "  +8.86%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "86.62"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +4.48%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "37.52"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +2.08%  "
$ws.Range("E26").Value = "  +8.94%  "
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "5.52"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -3.37%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "9.54"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +19.41%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "9.77"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +5.35%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "751.98"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +7.34%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "13.80"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +13.18%  "
$ws.Range("E32").Value = "  +12.66%  "
$ws.Range("E33").Value = "  -1.02%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "43.39"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +13.26%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.155"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +3.45%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "57.68"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +5.89%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.14%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "5.53"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +15.90%  "
$ws.Range("E39").Value = "  +6.16%  "
$ws.Range("B40").Value = "ThetaToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.90"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.88%  "
$ws.Range("D41").Value = "0.0₃0690"
$ws.Range("E41").Value = "  -9.49%  "
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.343"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +17.73%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.141"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +5.10%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.04%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "3.47"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +5.43%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "3.22"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +5.16%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "2.14"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +6.19%  "
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.45"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +12.63%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "146.75"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +1.52%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "2.67"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +5.21%  "
$ws.Range("E51").Value = "  +5.27%  "
